# Add the new "T_sup" column (E) with its header and data, matching the
# existing A:D layout (header row 1, data rows 2-9).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E1").Value = "T_sup"

$values = @(15.91, 13.64, 15.83, 15.55, 12.78, 13.66, 12.47, 11.42)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $values[$i]
}

# Mirror the author's final selection (whole new column's data range,
# anchored at the top cell) recorded in the saved sheet view.
$null = $ws.Range("E2:E9").Select()
